# My Courses Distribution Scenario for Available Courses Completed
#
# - The "Forgot Password" suite's Runmode flips from YES to NO.
# - A new test-suite row "MC Distribution" is appended (row 6).
# - Selection/active cell moves to B17 to match the post-edit workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Forgot Password row (row 4): Runmode YES -> NO
$ws.Range("C4").Value = "NO"

# New row 6: MC Distribution suite
$ws.Range("A6").Value = "MC Distribution"
$ws.Range("B6").Value = "My Courses distribution description"
$ws.Range("C6").Value = "Yes"

# Match the resulting selection recorded in the sheet view
$ws.Range("B17").Select()

# Best-effort window sizing to mirror the recorded workbook view (cosmetic)
$excel.ActiveWindow.Width = 13395
$excel.ActiveWindow.Height = 4215
